# merge jdm file into lda.R for testing
# - Update the B-column values for data rows 2-91 (LDA topic assignments)
# - Remove the now-unused trailing rows 92-101 (the corresponding "91".."100"
#   shared-string labels in column A are dropped automatically once those
#   rows are gone)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B2:B91 values, in row order.
$newValues = @(
    4.0, 4.0, 2.0, 4.0, 2.0, 2.0, 1.0, 2.0, 2.0, 5.0,
    5.0, 2.0, 5.0, 5.0, 4.0, 3.0, 3.0, 4.0, 4.0, 4.0,
    5.0, 3.0, 3.0, 5.0, 1.0, 2.0, 5.0, 1.0, 4.0, 5.0,
    1.0, 4.0, 3.0, 2.0, 4.0, 2.0, 2.0, 5.0, 3.0, 4.0,
    3.0, 4.0, 2.0, 4.0, 5.0, 3.0, 2.0, 3.0, 2.0, 5.0,
    3.0, 3.0, 4.0, 5.0, 5.0, 5.0, 1.0, 2.0, 5.0, 3.0,
    4.0, 3.0, 3.0, 2.0, 2.0, 5.0, 3.0, 2.0, 5.0, 2.0,
    2.0, 2.0, 2.0, 1.0, 5.0, 2.0, 4.0, 5.0, 5.0, 2.0,
    4.0, 4.0, 5.0, 1.0, 5.0, 3.0, 5.0, 5.0, 2.0, 5.0
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

# Drop the trailing rows (old rows 92-101), which shrinks sharedStrings /
# column A down to just "x","1",...,"90".
$ws.Range("A92:B101").EntireRow.Delete()
